$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.093.70"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "1.655.42"
$ws.Range("E3").Value = "  +3.42%  "
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'215.78"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "'0.508"
$ws.Range("E6").Value = "  +1.16%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("D10").Value = "'19.53"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("D11").Value = "'0.0859"
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("D12").Value = "1.885.53"
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "1.661.09"
$ws.Range("E13").Value = "  +3.58%  "
$ws.Range("D14").Value = "'4.07"
$ws.Range("E14").Value = "  +1.62%  "
$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = "  +2.66%  "
$ws.Range("D16").Value = "'65.09"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "'242.61"
$ws.Range("E17").Value = "  +5.69%  "
$ws.Range("D18").Value = "27.045.30"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").Value = "'7.85"
$ws.Range("E19").Value = "  +3.34%  "
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'4.43"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("D23").Value = "'2.29"
$ws.Range("E23").Value = "  +5.45%  "
$ws.Range("D24").Value = "'9.26"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'7.12"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("E28").Value = "  +1.09%  "
$ws.Range("D29").Value = "'15.88"
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("B32").Value = "Maker"
$ws.Range("C32").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D32").Value = "1.523.71"
$ws.Range("E32").Value = "  +5.42%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.29"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("D34").Value = "'3.05"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("D36").Value = "'2.43"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'0.578"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").Value = "'0.889"
$ws.Range("E38").Value = "  +7.87%  "
$ws.Range("D39").Value = "'0.0169"
$ws.Range("E39").Value = "  +2.33%  "
$ws.Range("D40").Value = "'5.99"
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("E42").Value = "  +2.56%  "
$ws.Range("D43").Value = "'64.89"
$ws.Range("E43").Value = "  +6.91%  "
$ws.Range("D44").Value = "1.793.22"
$ws.Range("E44").Value = "  +3.09%  "
$ws.Range("D45").Value = "'0.772"
$ws.Range("E45").Value = "  +1.83%  "
$ws.Range("D46").Value = "'0.911"
$ws.Range("E46").Value = "  -0.93%  "
$ws.Range("D47").Value = "'90.45"
$ws.Range("E47").Value = "  +3.48%  "
$ws.Range("E48").Value = "  +3.92%  "
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("D51").Value = "'7.52"
$ws.Range("E51").Value = "  +1.84%  "

# Restore default (General/Normal) style on cells that needed a text-force
# apostrophe prefix, so number formatting/style stays identical to the source.
$textForcedCells = @("D4","D5","D6","D7","D10","D11","D14","D15","D16","D17","D19","D22","D23","D24","D26","D27","D29","D33","D34","D35","D36","D37","D38","D39","D40","D41","D43","D45","D46","D47","D51")
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
